$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A27").Value = 1
$ws.Range("B27").Value = "VAWQ6"
$ws.Range("C27").Value = "VAWQ6"
$ws.Range("D27").Value = "VAWQ6"
$ws.Range("E27").Value = "U8"
$ws.Range("F27").Value = "CONV DC/DC DUAL 12V 230MA 6W"
$ws.Range("G27").Value = "Digikey"
$ws.Range("H27").Value = "102-1319-ND"
$ws.Range("I27").Value = "VAWQ6-Q24-D12H"
$ws.Range("J27").Value = "CUI Inc"
